$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H33").Value = 176.25
$ws.Range("I33").Value = 176.25
$ws.Range("K33").Value = 176.25
$ws.Range("M33").Value = 52.75
$ws.Range("H132").Value = 2743.6667
$ws.Range("I132").Value = 1190.625
$ws.Range("K132").Value = 3571.875
$ws.Range("M132").Value = -1041.875
$ws.Range("H135").Value = 636.9
$ws.Range("I135").Value = 636.9
$ws.Range("K135").Value = 5732.099999999999
$ws.Range("M135").Value = -3197.099999999999
$ws.Range("H137").Value = 2793.75
$ws.Range("J137").Value = 3694.818
$ws.Range("L137").Value = 11084.454
$ws.Range("N137").Value = -16184.454
$ws.Range("H141").Value = 8000
$ws.Range("I141").Value = 8000
$ws.Range("K141").Value = 24000
$ws.Range("M141").Value = -18820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4707.5
$ws.Range("I61").Value = 4707.5
$ws.Range("K61").Value = 4707.5
$ws.Range("M61").Value = -4495.5
$ws.Range("H136").Value = 4707.5
$ws.Range("I136").Value = 4707.5
$ws.Range("K136").Value = 14122.5
$ws.Range("M136").Value = -11572.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H80").Value = 819.25
$ws.Range("I80").Value = 53.666668
$ws.Range("K80").Value = 53.666668
$ws.Range("M80").Value = 944.333332
$ws.Range("H83").Value = 819.25
$ws.Range("I83").Value = 53.666668
$ws.Range("K83").Value = 268.33334
$ws.Range("M83").Value = 4723.66666
$ws.Range("H134").Value = 2791.8
$ws.Range("I134").Value = 2791.8
$ws.Range("K134").Value = 8375.400000000001
$ws.Range("M134").Value = -5840.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 5883
$ws.Range("J33").Value = 7735
$ws.Range("L33").Value = 7735
$ws.Range("N33").Value = -8493
$ws.Range("H39").Value = 5051
$ws.Range("I39").Value = 5051
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 5051
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -4660
$ws.Range("N39").ClearContents()
$ws.Range("H49").Value = 5051
$ws.Range("I49").Value = 5051
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 5051
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -4869
$ws.Range("N49").ClearContents()
$ws.Range("H50").Value = 53325
$ws.Range("J50").Value = 53325
$ws.Range("L50").Value = 53325
$ws.Range("N50").Value = -54575
$ws.Range("H58").Value = 3539.889
$ws.Range("J58").Value = 2889.6667
$ws.Range("L58").Value = 2889.6667
$ws.Range("N58").Value = -3295.6667
$ws.Range("H59").Value = 44418
$ws.Range("I59").Value = 5000
$ws.Range("J59").Value = 64127
$ws.Range("K59").Value = 5000
$ws.Range("L59").Value = 64127
$ws.Range("M59").Value = -3855
$ws.Range("N59").Value = -66417
$ws.Range("H86").Value = 27748.75
$ws.Range("I86").Value = 52500
$ws.Range("J86").Value = 2997.5
$ws.Range("K86").Value = 52500
$ws.Range("L86").Value = 2997.5
$ws.Range("M86").Value = -51377
$ws.Range("N86").Value = -5243.5
$ws.Range("H89").Value = 27748.75
$ws.Range("I89").Value = 52500
$ws.Range("J89").Value = 2997.5
$ws.Range("K89").Value = 262500
$ws.Range("L89").Value = 14987.5
$ws.Range("M89").Value = -256884
$ws.Range("N89").Value = -26219.5
$ws.Range("H136").Value = 3539.889
$ws.Range("J136").Value = 2889.6667
$ws.Range("L136").Value = 8669.000100000001
$ws.Range("N136").Value = -13769.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 2394.5715
$ws.Range("I34").Value = 442.4
$ws.Range("J34").Value = 3479.111
$ws.Range("K34").Value = 1327.2
$ws.Range("L34").Value = 10437.333
$ws.Range("M34").Value = -1243.2
$ws.Range("N34").Value = -10605.333
$ws.Range("H38").Value = 93.333336
$ws.Range("I38").Value = 93.333336
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 280.000008
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 66.99999200000002
$ws.Range("N38").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H40").Value = 188.5
$ws.Range("I40").Value = 100
$ws.Range("J40").Value = 218
$ws.Range("K40").Value = 400
$ws.Range("L40").Value = 872
$ws.Range("M40").Value = -331
$ws.Range("N40").Value = -1010
$ws.Range("H55").Value = 1862.25
$ws.Range("I55").Value = 1324.5
$ws.Range("J55").Value = 2400
$ws.Range("K55").Value = 3973.5
$ws.Range("L55").Value = 7200
$ws.Range("M55").Value = -3796.5
$ws.Range("N55").Value = -7554
$ws.Range("H68").Value = 812.5
$ws.Range("J68").Value = 833.3333
$ws.Range("L68").Value = 2499.9999
$ws.Range("N68").Value = -4121.9999
$ws.Range("H71").Value = 812.5
$ws.Range("J71").Value = 833.3333
$ws.Range("L71").Value = 7499.9997
$ws.Range("N71").Value = -15611.9997
$ws.Range("H82").Value = 14000
$ws.Range("J82").Value = 18000
$ws.Range("L82").Value = 54000
$ws.Range("N82").Value = -54812
$ws.Range("H85").Value = 14000
$ws.Range("J85").Value = 18000
$ws.Range("L85").Value = 54000
$ws.Range("N85").Value = -56808
$ws.Range("H104").Value = 5125.231
$ws.Range("J104").Value = 5170.24
$ws.Range("L104").Value = 15510.72
$ws.Range("N104").Value = -20752.72
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H117").Value = 1178.6666
$ws.Range("I117").Value = 420
$ws.Range("J117").Value = 1330.4
$ws.Range("K117").Value = 1260
$ws.Range("L117").Value = 3991.2
$ws.Range("M117").Value = 2182
$ws.Range("N117").Value = -10875.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2815.6667
$ws.Range("I97").Value = 1792.625
$ws.Range("K97").Value = 1792.625
$ws.Range("M97").Value = -1296.625
$ws.Range("H126").Value = 3993
$ws.Range("I126").Value = 3989.6667
$ws.Range("K126").Value = 11969.0001
$ws.Range("M126").Value = -9499.000100000001
$ws.Range("H132").Value = 3434.6
$ws.Range("I132").Value = 1808
$ws.Range("K132").Value = 5424
$ws.Range("M132").Value = -2894

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2758.75
$ws.Range("I7").Value = 3178.3333
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 3178.3333
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -3066.3333
$ws.Range("N7").Value = -1724
$ws.Range("H61").Value = 10202099
$ws.Range("I61").Value = 8502482
$ws.Range("K61").Value = 8502482
$ws.Range("M61").Value = -8502280
$ws.Range("H82").Value = 3877.4614
$ws.Range("I82").Value = 1544.8889
$ws.Range("J82").Value = 9125.75
$ws.Range("K82").Value = 1544.8889
$ws.Range("L82").Value = 9125.75
$ws.Range("M82").Value = -1183.8889
$ws.Range("N82").Value = -9847.75
$ws.Range("H85").Value = 3877.4614
$ws.Range("I85").Value = 1544.8889
$ws.Range("J85").Value = 9125.75
$ws.Range("K85").Value = 1544.8889
$ws.Range("L85").Value = 9125.75
$ws.Range("M85").Value = -296.8888999999999
$ws.Range("N85").Value = -11621.75
$ws.Range("H113").Value = 10202099
$ws.Range("I113").Value = 8502482
$ws.Range("K113").Value = 8502482
$ws.Range("M113").Value = -8500312
$ws.Range("H126").Value = 2758.75
$ws.Range("I126").Value = 3178.3333
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 9534.999899999999
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -7064.999899999999
$ws.Range("N126").Value = -9440
$ws.Range("H136").Value = 7643.625
$ws.Range("I136").Value = 3691.5
$ws.Range("K136").Value = 11074.5
$ws.Range("M136").Value = -8524.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 822.8
$ws.Range("I113").Value = 708.1429000000001
$ws.Range("J113").Value = 923.125
$ws.Range("K113").Value = 2124.4287
$ws.Range("L113").Value = 2769.375
$ws.Range("M113").Value = 45.57129999999961
$ws.Range("N113").Value = -7109.375
$ws.Range("H136").Value = 3553.5
$ws.Range("I136").Value = 3553.5
$ws.Range("K136").Value = 10660.5
$ws.Range("M136").Value = -8110.5
